$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to plain text so values like "1.011" or "0.06408"
# are not re-interpreted as numbers by Excels automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.878.32"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.643.53"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "216.37"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "0.5043"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "0.06408"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "19.77"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "0.07757"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "1.655.54"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "4.287"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "1.865.18"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "0.5469"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "0.0₅7941"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "64.08"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "25.901.63"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "201.88"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").Value = "4.374"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").Value = "9.956"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "5.991"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "1.928"
$ws.Range("E25").Value = "  +10.98%  "
$ws.Range("D26").Value = "142.12"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "0.1139"
$ws.Range("E27").Value = "  -3.57%  "
$ws.Range("D28").Value = "15.72"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "6.736"
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").Value = "1.248"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "0.05002"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").Value = "3.284"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").Value = "3.209"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "1.546"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "2.384"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "1.175.43"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "2.643"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("D38").Value = "0.8952"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").Value = "0.5599"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").Value = "0.01565"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "1.008"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "5.715"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D45").Value = "1.774.34"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").Value = "0.4552"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D50").Value = "0.05086"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").Value = "  -0.05%  "

# Row 43/44 swap (TrustWalletToken now ranked above Quant)
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8106"
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "99.95"
$ws.Range("E44").Value = "  -0.43%  "

# Row 48/49 swap (Aave now ranked above Frax)
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "55.14"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.16%  "
